# "Generate Report for Handoff" -- refresh the localization-status report:
#   - flip the per-language status from "Handed back: in sync with en-US"
#     to "Ready for handoff" (Overview!E2:F2 and the Status column on each
#     language sheet)
#   - bump the "generated" timestamps that go with that status
#   - the status text got much shorter, so the status columns shrink to fit

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

# ---- Overview sheet -------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("G2").Value = "2016-08-19 15:06:04"

# ---- zh-cn sheet ------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $newStatus
$zh.Range("H2").Value = "2016-08-19 15:05:56"

# ---- de-de sheet ------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $newStatus
$de.Range("H2").Value = "2016-08-19 15:06:04"

# ---- shrink the now-too-wide Status columns to fit the shorter text ---
# (Overview columns E & F, and column C on each language sheet)
$newStatusWidth = 16.33
$ov.Columns.Item(5).ColumnWidth = $newStatusWidth
$ov.Columns.Item(6).ColumnWidth = $newStatusWidth
$zh.Columns.Item(3).ColumnWidth = $newStatusWidth
$de.Columns.Item(3).ColumnWidth = $newStatusWidth
